$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 240
$ws.Range("D3").Value = 44839
$ws.Range("J3").Value = 200
$ws.Range("D4").Value = 44883
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 700
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = 750
$ws.Range("P4").Value = 750
$ws.Range("D5").Value = 44883
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 600
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = 600
$ws.Range("P5").Value = 600
$ws.Range("D6").Value = 44624
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 650
$ws.Range("L6").Value = 700
$ws.Range("M6").Value = 675
$ws.Range("P6").Value = 675
$ws.Range("D7").Value = 44882
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 750
$ws.Range("P7").Value = 750
$ws.Range("D8").Value = 44882
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = 600
$ws.Range("P8").Value = 600
$ws.Range("D9").Value = 44799
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 750
$ws.Range("L9").Value = 850
$ws.Range("M9").Value = 800
$ws.Range("P9").Value = 800
$ws.Range("D10").Value = 44799
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 650
$ws.Range("L10").Value = 650
$ws.Range("M10").Value = 650
$ws.Range("P10").Value = 650
$ws.Range("D11").Value = 44831
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 700
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = 750
$ws.Range("P11").Value = 750
$ws.Range("D12").Value = 44831
$ws.Range("I12").Value = "Segunda"
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 600
$ws.Range("M12").Value = 600
$ws.Range("P12").Value = 600
$ws.Range("D13").Value = 44608
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 600
$ws.Range("M13").Value = 625
$ws.Range("P13").Value = 625
$ws.Range("D16").Value = 44818
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = 850
$ws.Range("P16").Value = 850
$ws.Range("D17").Value = 44804
$ws.Range("I17").Value = "Primera"
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 850
$ws.Range("M17").Value = 800
$ws.Range("P17").Value = 800
$ws.Range("D18").Value = 44804
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 650
$ws.Range("L18").Value = 650
$ws.Range("M18").Value = 650
$ws.Range("P18").Value = 650
$ws.Range("D19").Value = 44859
$ws.Range("D20").Value = 44859
$ws.Range("D21").Value = 44764
$ws.Range("D22").Value = 44764
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 600
$ws.Range("P22").Value = 600
$ws.Range("D23").Value = 44791
$ws.Range("J23").Value = 240
$ws.Range("K23").Value = 750
$ws.Range("M23").Value = 775
$ws.Range("P23").Value = 775
$ws.Range("D24").Value = 44791
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 650
$ws.Range("L24").Value = 650
$ws.Range("M24").Value = 650
$ws.Range("P24").Value = 650
$ws.Range("D25").Value = 44797
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = 750
$ws.Range("L25").Value = 850
$ws.Range("M25").Value = 800
$ws.Range("P25").Value = 800
$ws.Range("D26").Value = 44797
$ws.Range("K26").Value = 650
$ws.Range("L26").Value = 650
$ws.Range("M26").Value = 650
$ws.Range("P26").Value = 650
$ws.Range("D27").Value = 44868
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = 750
$ws.Range("P27").Value = 750
$ws.Range("D28").Value = 44811
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 750
$ws.Range("L28").Value = 850
$ws.Range("M28").Value = 800
$ws.Range("P28").Value = 800
$ws.Range("D30").Value = 44754
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 700
$ws.Range("L30").Value = 750
$ws.Range("M30").Value = 725
$ws.Range("P30").Value = 725
$ws.Range("D31").Value = 44837
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 700
$ws.Range("L31").Value = 800
$ws.Range("M31").Value = 750
$ws.Range("P31").Value = 750
$ws.Range("D32").Value = 44837
$ws.Range("J32").Value = 150
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = 600
$ws.Range("P32").Value = 600
